# Terms and Conditions.docx edit script
# Applies:
#  1) Grammar-check style run-splitting (w:proofErr gramStart/gramEnd) around
#     a few words in several existing paragraphs (no visible text change).
#  2) Bolds "stored" inside the "Your data being processed..." paragraph and
#     "system" inside the "When using the system you will not..." paragraph
#     (already-bold paragraphs, split into runs the same way Word would).
#  3) Appends two new paragraphs (each preceded by a blank paragraph) that
#     contain a small rectangle "checkbox" shape followed by
#     "If you agree to adhere to these terms and conditions please click here:"
#     and "If you do not agree to these terms and conditions please click here:"

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$pkgFooter = '</pkg:xmlData></pkg:part></pkg:package>'

$docOpen = '<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"><w:body>'
$docClose = '</w:body></w:document>'

function Apply-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $full = $pkgHeader + $docOpen + $innerXml + $docClose + $pkgFooter
    $rng.InsertXML($full)
}

# ---------------------------------------------------------------------------
# Paragraph 3: "Please read the following terms and conditions carefully..."
# Split off "America" with proofErr marks.
# ---------------------------------------------------------------------------
$para3 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Please read the following terms and conditions carefully before progressing onto the system. These terms will be applicable to the law of the following countries/regions, North </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>America</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> and Europe. </w:t>
  </w:r>
</w:p>
'@
Apply-ParaXml 3 $para3

# ---------------------------------------------------------------------------
# Paragraph 4: "In order to access the system, you must first agree..."
# Split off "In order to" and "particular group" with proofErr marks.
# ---------------------------------------------------------------------------
$para4 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>In order to</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> access </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">the system, you must first agree to a range of conditions set out by the ownership of the system. Any violations of the terms and conditions can result in a prosecution such as a large fine or even a jail sentence, depending on the laws of the country. These terms apply to all users of the system and are not limited to any </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>particular group</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> or individual.</w:t>
  </w:r>
</w:p>
'@
Apply-ParaXml 4 $para4

# ---------------------------------------------------------------------------
# Paragraph 5: "By accessing this system you agree to respect..." (bold)
# Split off "system" with proofErr marks.
# ---------------------------------------------------------------------------
$para5 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">By accessing this </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>system</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> you agree to respect these terms no matter the circumstance. However, if there are any terms which you disagree with, you may not access the system.</w:t>
  </w:r>
</w:p>
'@
Apply-ParaXml 5 $para5

# ---------------------------------------------------------------------------
# Paragraph 7: "Your data being processed, stored and sent..." (bold)
# Split off "stored" with proofErr marks.
# ---------------------------------------------------------------------------
$para7 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Your data being processed, </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>stored</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> and sent for marketing purposes-</w:t>
  </w:r>
</w:p>
'@
Apply-ParaXml 7 $para7

# ---------------------------------------------------------------------------
# Paragraph 11: "When using the system you will not be permitted..." (bold)
# Split off "system" with proofErr marks.
# ---------------------------------------------------------------------------
$para11 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">When using the </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>system</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> you will not be permitted to modify or alter any of the features or functionalities of the system-</w:t>
  </w:r>
</w:p>
'@
Apply-ParaXml 11 $para11

# ---------------------------------------------------------------------------
# Append two new paragraphs (each with a leading blank paragraph) after the
# last paragraph of body text ("When using the system, you will not be
# allowed to make any changes ..."). Each new paragraph contains a small
# rectangle "checkbox" drawing followed by click-here text.
# ---------------------------------------------------------------------------
$newBlock = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:noProof/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <mc:AlternateContent>
      <mc:Choice Requires="wps">
        <w:drawing>
          <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="2C4810AE" wp14:editId="08EFB2BD">
            <wp:simplePos x="0" y="0"/>
            <wp:positionH relativeFrom="margin">
              <wp:align>left</wp:align>
            </wp:positionH>
            <wp:positionV relativeFrom="paragraph">
              <wp:posOffset>290195</wp:posOffset>
            </wp:positionV>
            <wp:extent cx="266700" cy="190500"/>
            <wp:effectExtent l="0" t="0" r="19050" b="19050"/>
            <wp:wrapNone/>
            <wp:docPr id="1" name="Rectangle 1"/>
            <wp:cNvGraphicFramePr/>
            <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
              <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                <wps:wsp>
                  <wps:cNvSpPr/>
                  <wps:spPr>
                    <a:xfrm>
                      <a:off x="0" y="0"/>
                      <a:ext cx="266700" cy="190500"/>
                    </a:xfrm>
                    <a:prstGeom prst="rect">
                      <a:avLst/>
                    </a:prstGeom>
                    <a:ln>
                      <a:solidFill>
                        <a:schemeClr val="tx1"/>
                      </a:solidFill>
                    </a:ln>
                  </wps:spPr>
                  <wps:style>
                    <a:lnRef idx="2">
                      <a:schemeClr val="accent6"/>
                    </a:lnRef>
                    <a:fillRef idx="1">
                      <a:schemeClr val="lt1"/>
                    </a:fillRef>
                    <a:effectRef idx="0">
                      <a:schemeClr val="accent6"/>
                    </a:effectRef>
                    <a:fontRef idx="minor">
                      <a:schemeClr val="dk1"/>
                    </a:fontRef>
                  </wps:style>
                  <wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1">
                    <a:prstTxWarp prst="textNoShape">
                      <a:avLst/>
                    </a:prstTxWarp>
                    <a:noAutofit/>
                  </wps:bodyPr>
                </wps:wsp>
              </a:graphicData>
            </a:graphic>
          </wp:anchor>
        </w:drawing>
      </mc:Choice>
      <mc:Fallback>
        <w:pict>
          <v:rect w14:anchorId="3A791F47" id="Rectangle 1" o:spid="_x0000_s1026" style="position:absolute;margin-left:0;margin-top:22.85pt;width:21pt;height:15pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:left;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQAU7AV/eAIAAFEFAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtv2zAMvg/YfxB0X+0EbboGdYqgRYcB&#10;RVv0gZ5VWUqEyaJGKXGyXz9Kdpysy2nYRSbN98fH5dWmsWytMBhwFR+dlJwpJ6E2blHx15fbL185&#10;C1G4WlhwquJbFfjV7POny9ZP1RiWYGuFjJy4MG19xZcx+mlRBLlUjQgn4JUjoQZsRCQWF0WNoiXv&#10;jS3GZTkpWsDaI0gVAv296YR8lv1rrWR80DqoyGzFKbeYX8zve3qL2aWYLlD4pZF9GuIfsmiEcRR0&#10;cHUjomArNH+5aoxECKDjiYSmAK2NVLkGqmZUfqjmeSm8yrUQOMEPMIX/51berx+RmZp6x5kTDbXo&#10;iUATbmEVGyV4Wh+mpPXsH7HnApGp1o3GJn2pCrbJkG4HSNUmMkk/x5PJeUnASxKNLsozoslLsTf2&#10;GOI3BQ1LRMWRgmcgxfouxE51p5JiWZfeANbUt8bazKRZUdcW2VpQl+Mmp00hDrSIS5ZFKqZLP1Nx&#10;a1Xn9UlpQiElnKPn+dv7FFIqFyd96taRdjLTlMFgODpmaOMumV43mak8l4Nheczwz4iDRY4KLg7G&#10;jXGAxxzUP4bInf6u+q7mVP471FtqPkK3FcHLW0NNuBMhPgqkNaC+0WrHB3q0hbbi0FOcLQF/Hfuf&#10;9Gk6ScpZS2tV8fBzJVBxZr87mtuL0elp2sPMnJ6dj4nBQ8n7ocStmmugntJsUnaZTPrR7kiN0LzR&#10;BZinqCQSTlLsisuIO+Y6dutON0Sq+Tyr0e55Ee/cs5fJeUI1DdnL5k2g7ycx0gjfw24FxfTDQHa6&#10;ydLBfBVBmzyte1x7vGlv87z3NyYdhkM+a+0v4ew3AAAA//8DAFBLAwQUAAYACAAAACEAeBcUGdoA&#10;AAAFAQAADwAAAGRycy9kb3ducmV2LnhtbEyPy07DMBBF90j8gzVI7KhDVSgKmVQVohJiASLlA9x4&#10;GkeNH9hOm/49w4ouj+7o3jPVarKDOFJMvXcI97MCBLnW6951CN/bzd0TiJSV02rwjhDOlGBVX19V&#10;qtT+5L7o2OROcIlLpUIwOYdSytQasirNfCDH2d5HqzJj7KSO6sTldpDzoriUVvWOF4wK9GKoPTSj&#10;RQhxHT7Nq9lupo/49t6NTW9+zoi3N9P6GUSmKf8fw58+q0PNTjs/Op3EgMCPZITFwxIEp4s58w5h&#10;ySzrSl7a178AAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAA&#10;AAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAA&#10;AAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAFOwFf3gCAABRBQAADgAAAAAA&#10;AAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAeBcUGdoAAAAFAQAADwAA&#10;AAAAAAAAAAAAAADSBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAANkFAAAAAA==&#10;" fillcolor="white [3201]" strokecolor="black [3213]" strokeweight="1pt">
            <w10:wrap anchorx="margin"/>
          </v:rect>
        </w:pict>
      </mc:Fallback>
    </mc:AlternateContent>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">If you agree to adhere to these terms and </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>conditions</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> please click here:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:noProof/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <mc:AlternateContent>
      <mc:Choice Requires="wps">
        <w:drawing>
          <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="54535831" wp14:editId="1A06EE65">
            <wp:simplePos x="0" y="0"/>
            <wp:positionH relativeFrom="margin">
              <wp:align>left</wp:align>
            </wp:positionH>
            <wp:positionV relativeFrom="paragraph">
              <wp:posOffset>280670</wp:posOffset>
            </wp:positionV>
            <wp:extent cx="266700" cy="190500"/>
            <wp:effectExtent l="0" t="0" r="19050" b="19050"/>
            <wp:wrapNone/>
            <wp:docPr id="2" name="Rectangle 2"/>
            <wp:cNvGraphicFramePr/>
            <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
              <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                <wps:wsp>
                  <wps:cNvSpPr/>
                  <wps:spPr>
                    <a:xfrm>
                      <a:off x="0" y="0"/>
                      <a:ext cx="266700" cy="190500"/>
                    </a:xfrm>
                    <a:prstGeom prst="rect">
                      <a:avLst/>
                    </a:prstGeom>
                    <a:ln>
                      <a:solidFill>
                        <a:schemeClr val="tx1"/>
                      </a:solidFill>
                    </a:ln>
                  </wps:spPr>
                  <wps:style>
                    <a:lnRef idx="2">
                      <a:schemeClr val="accent6"/>
                    </a:lnRef>
                    <a:fillRef idx="1">
                      <a:schemeClr val="lt1"/>
                    </a:fillRef>
                    <a:effectRef idx="0">
                      <a:schemeClr val="accent6"/>
                    </a:effectRef>
                    <a:fontRef idx="minor">
                      <a:schemeClr val="dk1"/>
                    </a:fontRef>
                  </wps:style>
                  <wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1">
                    <a:prstTxWarp prst="textNoShape">
                      <a:avLst/>
                    </a:prstTxWarp>
                    <a:noAutofit/>
                  </wps:bodyPr>
                </wps:wsp>
              </a:graphicData>
            </a:graphic>
          </wp:anchor>
        </w:drawing>
      </mc:Choice>
      <mc:Fallback>
        <w:pict>
          <v:rect w14:anchorId="6BB1124B" id="Rectangle 2" o:spid="_x0000_s1026" style="position:absolute;margin-left:0;margin-top:22.1pt;width:21pt;height:15pt;z-index:251661312;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:left;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDn6hzKeAIAAFEFAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtv2zAMvg/YfxB0X20HbboGdYogRYcB&#10;RVv0gZ4VWUqESaImKXGyXz9Kdpysy2nYRSZN8uOb1zdbo8lG+KDA1rQ6KykRlkOj7LKmb693X75S&#10;EiKzDdNgRU13ItCb6edP162biBGsQDfCEwSxYdK6mq5idJOiCHwlDAtn4IRFoQRvWETWL4vGsxbR&#10;jS5GZTkuWvCN88BFCPj3thPSacaXUvD4KGUQkeiaYmwxvz6/i/QW02s2WXrmVor3YbB/iMIwZdHp&#10;AHXLIiNrr/6CMop7CCDjGQdTgJSKi5wDZlOVH7J5WTEnci5YnOCGMoX/B8sfNk+eqKamI0osM9ii&#10;Zywas0styCiVp3Vhglov7sn3XEAy5bqV3qQvZkG2uaS7oaRiGwnHn6Px+LLEwnMUVVflBdKIUhyM&#10;nQ/xmwBDElFTj85zIdnmPsROda+SfGmb3gBaNXdK68ykWRFz7cmGYZfjtupdHGmhw2RZpGS68DMV&#10;d1p0qM9CYhVSwNl7nr8DJuNc2DjucbVF7WQmMYLBsDplqOM+mF43mYk8l4NhecrwT4+DRfYKNg7G&#10;RlnwpwCaH4PnTn+ffZdzSn8BzQ6b76HbiuD4ncIm3LMQn5jHNcC+4WrHR3ykhram0FOUrMD/OvU/&#10;6eN0opSSFteqpuHnmnlBif5ucW6vqvPztIeZOb+4HCHjjyWLY4ldmzlgTys8Io5nMulHvSelB/OO&#10;F2CWvKKIWY6+a8qj3zPz2K073hAuZrOshrvnWLy3L44n8FTVNGSv23fmXT+JEUf4AfYryCYfBrLT&#10;TZYWZusIUuVpPdS1rzfubZ73/sakw3DMZ63DJZz+BgAA//8DAFBLAwQUAAYACAAAACEAezGNkNkA&#10;AAAFAQAADwAAAGRycy9kb3ducmV2LnhtbEyP3UrDQBCF7wXfYRnBO7sxBC0xk1LEgnihmPoA2+yY&#10;DWZ/3N206ds7XunlxxnO+abZLHYSR4pp9A7hdlWAINd7PboB4WO/u1mDSFk5rSbvCOFMCTbt5UWj&#10;au1P7p2OXR4El7hUKwSTc6ilTL0hq9LKB3KcffpoVWaMg9RRnbjcTrIsijtp1eh4wahAj4b6r262&#10;CCFuw5t5Mvvd8hqfX4a5G833GfH6atk+gMi05L9j+NVndWjZ6eBnp5OYEPiRjFBVJQhOq5L5gHDP&#10;LNtG/rdvfwAAAP//AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAA&#10;AAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAA&#10;AAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQDn6hzKeAIAAFEFAAAOAAAAAAAA&#10;AAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQB7MY2Q2QAAAAUBAAAPAAAA&#10;AAAAAAAAAAAAANIEAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAA2AUAAAAA&#10;" fillcolor="white [3201]" strokecolor="black [3213]" strokeweight="1pt">
            <w10:wrap anchorx="margin"/>
          </v:rect>
        </w:pict>
      </mc:Fallback>
    </mc:AlternateContent>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">If you do not agree to these terms and </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>conditions</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> please click here:</w:t>
  </w:r>
</w:p>
'@

# Paragraph 12 is "When using the system, you will not be allowed to make
# any changes to the actual system, if this is violated, access to the
# system may be revoked or restricted. " -- the last paragraph of body text
# before the trailing blank paragraphs. Paragraph 13 is the first of those
# (originally empty) trailing paragraphs. Collapse *its* range to the start
# and insert the new block there, so paragraph 12's own content is left
# untouched and the new paragraphs land right after it.
$p13 = $d.Paragraphs.Item(13)
$insertPoint = $p13.Range
$insertPoint.Collapse(1)
$fullBlock = $pkgHeader + $docOpen + $newBlock + $docClose + $pkgFooter
$insertPoint.InsertXML($fullBlock)

Write-Output "Edit complete"
